$d = $word.ActiveDocument

function Get-ParaByText($text) {
    $r = $d.Content
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r.Paragraphs(1)
}

# ---------------------------------------------------------------------------
# 1. "Charlie Frey" heading -> bold
# ---------------------------------------------------------------------------
(Get-ParaByText("Charlie Frey")).Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2. Address line: drop the street address, keep "Colorado Springs, CO 80951"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("1945 Riverwalk Pkwy, Colorado Springs, CO 80951", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Colorado Springs, CO 80951", 2)

# ---------------------------------------------------------------------------
# 3. Contact block reorder: charliefrey.io <-> linkedin.com/in/csfrey95,
#    with a new "github.com/csfrey" line inserted between them. This part of
#    the document is always the first handful of paragraphs, so address the
#    paragraphs positionally (the text alone is no longer unique once we
#    start swapping it around).
# ---------------------------------------------------------------------------
$d.Paragraphs(5).Range.Text = "linkedin.com/in/csfrey95"
$d.Paragraphs(6).Range.Text = "charliefrey.io"
$d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.Text = "github.com/csfrey"

# ---------------------------------------------------------------------------
# 4. "Professional Summary" heading -> bold
# ---------------------------------------------------------------------------
(Get-ParaByText("Professional Summary")).Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 5. "Work Experience" heading -> bold
# ---------------------------------------------------------------------------
(Get-ParaByText("Work Experience")).Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 6. New DataAnnotation job block, inserted right before "Career Gap"
# ---------------------------------------------------------------------------
$anchor = Get-ParaByText("Career Gap")
$startIdx = $anchor.Index

for ($i = 0; $i -lt 7; $i++) {
    $anchor.Range.InsertParagraphBefore()
}

$lines = @(
    "DataAnnotation",
    "Data Annotation Specialist",
    "Remote | 07/2023 – Present",
    "- Conduct high-quality data labeling to train machine learning models, contributing to the improvement of AI systems.",
    "- Select and manage diverse projects, ensuring completion within tight deadlines while maintaining high accuracy and quality standards.",
    "- Effectively manage a self-directed work schedule, balancing multiple projects and personal commitments to consistently meet deadlines.",
    ""
)
for ($i = 0; $i -lt $lines.Length; $i++) {
    $d.Paragraphs($startIdx + $i).Range.Text = $lines[$i]
}

# ---------------------------------------------------------------------------
# 7. Date range update for the "Career Gap" / Software Engineer entry
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Colorado Springs, CO | 06/2022 - Present", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Colorado Springs, CO | 06/2022 - 07/2024", 2)

# ---------------------------------------------------------------------------
# 8. "Education" heading -> bold
# ---------------------------------------------------------------------------
(Get-ParaByText("Education")).Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 9. "Projects" heading -> bold
# ---------------------------------------------------------------------------
(Get-ParaByText("Projects")).Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 10. "Skills" heading -> bold
# ---------------------------------------------------------------------------
(Get-ParaByText("Skills")).Range.Font.Bold = 1

Write-Output "edit complete"
